$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$Bvals = @(1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02)
$Cvals = @(1.036308532144776,1.03714272496299,1.037683207584575,1.037910593365402,1.037948782176296,1.037686245270846,1.036590304720372,1.034664584035884,1.033384545539705,1.032831191225701,1.032625789324809,1.032669842433853,1.032814209804741,1.032903177713834,1.033421289176963,1.033746531850018,1.033936328164049,1.034001058637496,1.033711627299064,1.032771693374581,1.032181521939218,1.032494306366085,1.03372739889222,1.035161770896792)
$Dvals = @(1.045447246956343,1.046110420238542,1.046540159482722,1.046720968900327,1.04675133616674,1.046542574891529,1.045671239634903,1.04414067808227,1.043123668111898,1.042684112537701,1.042520966554854,1.042555956256867,1.042670624288736,1.042741291604717,1.043152857328752,1.043411241539477,1.043562031329046,1.043613460034319,1.043383511210296,1.042636853958802,1.042168122316547,1.042416536825271,1.043396041109643,1.044535779877572)
$Evals = @(1.039956798640218,1.040739082067003,1.041246306483762,1.041459789455829,1.041495648546712,1.041249158090121,1.040220959603455,1.038417150873809,1.03722010897776,1.036703104376612,1.036511266469538,1.036552407242131,1.036687242891536,1.036770346227573,1.037254448898294,1.037558469137878,1.037735926466139,1.03779645639061,1.037525837444979,1.036647531599766,1.036096467548631,1.036388486265836,1.03754058190879,1.038882518002286)
$Fvals = @(1.053993632004908,1.054787798995624,1.055302645883748,1.05551931721606,1.055555710698378,1.055305540156612,1.054261822891036,1.052430158988653,1.051214210053968,1.050688939725317,1.050494019804498,1.050535822209932,1.050672823719277,1.050757259958285,1.051249096890134,1.051557947581026,1.051738214815243,1.051799701550651,1.05152479841336,1.050632474959477,1.050072529519137,1.05036926283996,1.051539776725285,1.052902787124777)
$Ivals = @(1.040818700314365,1.041018470270297,1.041146931474414,1.041200743509382,1.041209767444516,1.041147651272955,1.040886379392688,1.040419865317077,1.04010478572192,1.039967396727032,1.039916221232702,1.039927205013798,1.039963169459694,1.039985309401293,1.04011388366791,1.040194279094643,1.040241079956875,1.040257022141185,1.040185662981332,1.039952582771182,1.039805208403042,1.039883412541632,1.040189556517748,1.040541191041227)
$Jvals = @(1.04141755587214,1.041896240001939,1.042206016512295,1.042336253866908,1.042358121683331,1.042207756723079,1.041579321219586,1.040472268823601,1.039734532908018,1.03941517214721,1.039296561222713,1.039322003019637,1.039405367440502,1.039456732889845,1.039755729705391,1.039943305848971,1.040052723834443,1.0400900338564,1.039923179867703,1.039380818302735,1.039039895018134,1.039220616794425,1.039932273907789,1.041018470270297)
$Kvals = @(1.048215461422746,1.048690463194736,1.048997758204943,1.049126928618031,1.049148615897251,1.048999484252495,1.048376002706641,1.047276926618737,1.046544006258202,1.046226612568155,1.046108714635509,1.046134004291444,1.046216867157509,1.046267921243789,1.046565070025976,1.046751455383616,1.046860167404068,1.046897234774266,1.046731458340983,1.046192466202178,1.045853559238437,1.046033221783939,1.046740494154484,1.047561105819007)
$Lvals = @(1.042740513027197,1.043333195163638,1.043717075810325,1.043878547651768,1.043905664625145,1.043719233059457,1.042940734031008,1.041571863703833,1.040661353903142,1.040267603362574,1.040121424483208,1.040152776812444,1.040255518578757,1.04031883150144,1.04068749662738,1.040918887142609,1.041053902053402,1.041099946868538,1.040894056077791,1.040225261534668,1.039805213616603,1.040027845816079,1.040905276016749,1.041925391215762)
$Mvals = @(1.056738067710037,1.057345438288602,1.057738756647028,1.057904180034992,1.057931959545267,1.05774096676018,1.056943266489945,1.0555400548145,1.054606316447395,1.054202429675768,1.054052474118898,1.054084637082875,1.054190032947251,1.054254979609402,1.054633130243209,1.054870449929168,1.055008915681561,1.055056135861824,1.054844983499935,1.054158994657637,1.053728068765673,1.05395647392307,1.054856490556045,1.055902519761575)
$Nvals = @(1.017742708499226,1.017903295395854,1.018007171422877,1.018050832181955,1.01805816249008,1.018007754854952,1.017796986434323,1.017425343022212,1.017177446289905,1.017070078418618,1.017030193648566,1.017038749220512,1.017066781600935,1.017084052818243,1.017184571415668,1.017247617204587,1.01728438811366,1.017296925581997,1.017240853256921,1.017058526858012,1.01694387064609,1.017004653841582,1.017243909601918,1.017521447544976)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $Bvals[$i]
    $ws.Cells.Item($row, 3).Value = $Cvals[$i]
    $ws.Cells.Item($row, 4).Value = $Dvals[$i]
    $ws.Cells.Item($row, 5).Value = $Evals[$i]
    $ws.Cells.Item($row, 6).Value = $Fvals[$i]
    $ws.Cells.Item($row, 9).Value = $Ivals[$i]
    $ws.Cells.Item($row, 10).Value = $Jvals[$i]
    $ws.Cells.Item($row, 11).Value = $Kvals[$i]
    $ws.Cells.Item($row, 12).Value = $Lvals[$i]
    $ws.Cells.Item($row, 13).Value = $Mvals[$i]
    $ws.Cells.Item($row, 14).Value = $Nvals[$i]
}